$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4246.2666
$ws.Range("I51").Value = 2167.3333
$ws.Range("K51").Value = 2167.3333
$ws.Range("M51").Value = -1683.3333

$ws.Range("H62").Value = 1853.5454
$ws.Range("I62").Value = 1738.9
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 1738.9
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -1114.9
$ws.Range("N62").Value = -4248

$ws.Range("H65").Value = 1853.5454
$ws.Range("I65").Value = 1738.9
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 8694.5
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -5574.5
$ws.Range("N65").Value = -21240

$ws.Range("H74").Value = 5775
$ws.Range("I74").Value = 3900
$ws.Range("K74").Value = 3900
$ws.Range("M74").Value = -2964

$ws.Range("H77").Value = 5775
$ws.Range("I77").Value = 3900
$ws.Range("K77").Value = 19500
$ws.Range("M77").Value = -14820

$ws.Range("H115").Value = 777
$ws.Range("I115").Value = 777
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 2331
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -764
$ws.Range("N115").ClearContents()

$ws.Range("H137").Value = 4711.643
$ws.Range("I137").Value = 2796.4827
$ws.Range("J137").Value = 8983.923000000001
$ws.Range("K137").Value = 8389.4481
$ws.Range("L137").Value = 26951.769
$ws.Range("M137").Value = -5839.4481
$ws.Range("N137").Value = -32051.769

$ws.Range("H141").Value = 1734.4546
$ws.Range("I141").Value = 1430.0555
$ws.Range("J141").Value = 3104.25
$ws.Range("K141").Value = 4290.166499999999
$ws.Range("L141").Value = 9312.75
$ws.Range("M141").Value = 889.8335000000006
$ws.Range("N141").Value = -19672.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3640.8208
$ws.Range("I32").Value = 3640.8208
$ws.Range("K32").Value = 3640.8208
$ws.Range("M32").Value = -3353.8208

$ws.Range("H74").Value = 52950.582
$ws.Range("I74").Value = 55187.566
$ws.Range("K74").Value = 55187.566
$ws.Range("M74").Value = -54313.566

$ws.Range("H77").Value = 52950.582
$ws.Range("I77").Value = 55187.566
$ws.Range("K77").Value = 275937.83
$ws.Range("M77").Value = -271569.83

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H122").Value = 1838.5946
$ws.Range("I122").Value = 1730.8064
$ws.Range("J122").Value = 2395.5
$ws.Range("K122").Value = 5192.4192
$ws.Range("L122").Value = 7186.5
$ws.Range("M122").Value = -2742.4192
$ws.Range("N122").Value = -12086.5

$ws.Range("H138").Value = 198999.5
$ws.Range("J138").Value = 198999.5
$ws.Range("L138").Value = 198999.5
$ws.Range("N138").Value = -209279.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2580.6
$ws.Range("I20").Value = 2302.5
$ws.Range("K20").Value = 2302.5
$ws.Range("M20").Value = -2055.5

$ws.Range("H99").Value = 3909.818
$ws.Range("I99").Value = 1999.8
$ws.Range("K99").Value = 1999.8
$ws.Range("M99").Value = -501.8

$ws.Range("H107").Value = 4700.4165
$ws.Range("I107").Value = 4493.6
$ws.Range("J107").Value = 5734.5
$ws.Range("K107").Value = 4493.6
$ws.Range("L107").Value = 5734.5
$ws.Range("M107").Value = -2573.6
$ws.Range("N107").Value = -9574.5

$ws.Range("H134").Value = 880.0454999999999
$ws.Range("I134").Value = 767.35
$ws.Range("K134").Value = 2302.05
$ws.Range("M134").Value = 232.9499999999998

$ws.Range("H137").Value = 134665.67
$ws.Range("J137").Value = 144750
$ws.Range("L137").Value = 144750
$ws.Range("N137").Value = -154950

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4075.4407
$ws.Range("I31").Value = 2463
$ws.Range("J31").Value = 6599.2607
$ws.Range("K31").Value = 2463
$ws.Range("L31").Value = 6599.2607
$ws.Range("M31").Value = -2168
$ws.Range("N31").Value = -7189.2607

$ws.Range("H34").Value = 4075.4407
$ws.Range("I34").Value = 2463
$ws.Range("J34").Value = 6599.2607
$ws.Range("K34").Value = 2463
$ws.Range("L34").Value = 6599.2607
$ws.Range("M34").Value = -2261
$ws.Range("N34").Value = -7003.2607

$ws.Range("H41").Value = 29909.6
$ws.Range("I41").Value = 12264.75
$ws.Range("J41").Value = 36325.91
$ws.Range("K41").Value = 12264.75
$ws.Range("L41").Value = 36325.91
$ws.Range("M41").Value = -11836.75
$ws.Range("N41").Value = -37181.91

$ws.Range("H86").Value = 8808.546
$ws.Range("I86").Value = 7049.6665
$ws.Range("K86").Value = 7049.6665
$ws.Range("M86").Value = -5926.6665

$ws.Range("H89").Value = 8808.546
$ws.Range("I89").Value = 7049.6665
$ws.Range("K89").Value = 35248.3325
$ws.Range("M89").Value = -29632.3325

$ws.Range("H112").Value = 60319
$ws.Range("I112").Value = 60319
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 60319
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -58842
$ws.Range("N112").ClearContents()

$ws.Range("H134").Value = 3007.4167
$ws.Range("I134").Value = 2548.875
$ws.Range("J134").Value = 3924.5
$ws.Range("K134").Value = 7646.625
$ws.Range("L134").Value = 11773.5
$ws.Range("M134").Value = -5111.625
$ws.Range("N134").Value = -16843.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 557521.9
$ws.Range("I32").Value = 502799.84
$ws.Range("J32").Value = 666966
$ws.Range("K32").Value = 1508399.52
$ws.Range("L32").Value = 2000898
$ws.Range("M32").Value = -1508116.52
$ws.Range("N32").Value = -2001464

$ws.Range("H106").Value = 4750.5
$ws.Range("J106").Value = 4750.5
$ws.Range("L106").Value = 14251.5
$ws.Range("N106").Value = -16143.5

$ws.Range("H129").Value = 1691.68
$ws.Range("J129").Value = 1695.591
$ws.Range("L129").Value = 5086.772999999999
$ws.Range("N129").Value = -15086.773

$ws.Range("H131").Value = 1508.909
$ws.Range("I131").Value = 1274.1666
$ws.Range("J131").Value = 1596.9375
$ws.Range("K131").Value = 3822.4998
$ws.Range("L131").Value = 4790.8125
$ws.Range("M131").Value = 1217.5002
$ws.Range("N131").Value = -14870.8125

$ws.Range("H132").Value = 2617.6875
$ws.Range("I132").Value = 1089.5
$ws.Range("J132").Value = 5164.6665
$ws.Range("K132").Value = 9805.5
$ws.Range("L132").Value = 46481.9985
$ws.Range("M132").Value = -7275.5
$ws.Range("N132").Value = -51541.9985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 1548009.5
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 1548009.5
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 1548009.5
$ws.Range("N33").Value = -1548513.5
$ws.Range("M33").ClearContents()

$ws.Range("H134").Value = 89969
$ws.Range("J134").Value = 89969
$ws.Range("L134").Value = 269907
$ws.Range("N134").Value = -274977

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 71432536
$ws.Range("I61").Value = 71432536
$ws.Range("K61").Value = 71432536
$ws.Range("M61").Value = -71432334

$ws.Range("H93").Value = 2379.8
$ws.Range("I93").Value = 2509.5557
$ws.Range("K93").Value = 2509.5557
$ws.Range("M93").Value = -1261.5557

$ws.Range("H113").Value = 71432536
$ws.Range("I113").Value = 71432536
$ws.Range("K113").Value = 71432536
$ws.Range("M113").Value = -71430366

$ws.Range("H119").Value = 74990
$ws.Range("J119").Value = 74990
$ws.Range("L119").Value = 74990
$ws.Range("N119").Value = -84666

$ws.Range("H121").Value = 94990
$ws.Range("J121").Value = 94990
$ws.Range("L121").Value = 94990
$ws.Range("N121").Value = -98484

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 132
$ws.Range("I14").Value = 132
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 132
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 36
$ws.Range("N14").ClearContents()

$ws.Range("H119").Value = 62485
$ws.Range("J119").Value = 62485
$ws.Range("L119").Value = 62485
$ws.Range("N119").Value = -72161

$ws.Range("H122").Value = 2175.44
$ws.Range("I122").Value = 2021
$ws.Range("J122").Value = 2986.25
$ws.Range("K122").Value = 6063
$ws.Range("L122").Value = 8958.75
$ws.Range("M122").Value = -3613
$ws.Range("N122").Value = -13858.75

$ws.Range("H136").Value = 4396.8423
$ws.Range("I136").Value = 2628.2917
$ws.Range("J136").Value = 13829.111
$ws.Range("K136").Value = 7884.875100000001
$ws.Range("L136").Value = 41487.333
$ws.Range("M136").Value = -5334.875100000001
$ws.Range("N136").Value = -46587.333
